# Generate Report for Handoff
#
# Regenerates the localization-status report for a new source file
# (24659008-6d08-43ac-8072-5cfee2563ff9.md), replacing every reference to
# the previous file (aecc829d-3045-4f88-8f95-3781bc097e39.md) across the
# "Overview", "zh-cn" and "de-de" sheets, and refreshing the handoff /
# generate timestamps that were recorded for the new run.

$wb = $excel.ActiveWorkbook

$oldGuid = "aecc829d-3045-4f88-8f95-3781bc097e39"
$newGuid = "24659008-6d08-43ac-8072-5cfee2563ff9"

$repoBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/29eb39c3e860fdc474da7fad1ab4db34fef85cbf/e2e/"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2016-08-20 01:00:55"

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("B2"),
    "$repoBase$newGuid.md",
    [System.Reflection.Missing]::Value,
    [System.Reflection.Missing]::Value,
    "e2e\$newGuid.md"
)

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = "$newGuid.md"
$wsZhCn.Range("G2").Value = "$newGuid.013fa165bbd107f6d1205b71446064b63bcc4385.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-20 01:00:51"

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A2"),
    "$repoBase$newGuid.md",
    [System.Reflection.Missing]::Value,
    [System.Reflection.Missing]::Value,
    "$newGuid.md"
)

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = "$newGuid.md"
$wsDeDe.Range("G2").Value = "$newGuid.013fa165bbd107f6d1205b71446064b63bcc4385.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-20 01:00:55"

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A2"),
    "$repoBase$newGuid.md",
    [System.Reflection.Missing]::Value,
    [System.Reflection.Missing]::Value,
    "$newGuid.md"
)

Write-Output "Localization status report regenerated for $newGuid.md"
